# Add two new "FAIL" rows (35 and 36) to the daily-update sheet,
# mirroring the layout/styling of the existing rows above (row 34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: 864954_GJ_P47 ---
$ws.Range("A35").Value = "GJ"
$ws.Range("B35").Value = "864954_GJ_P47"
$ws.Range("C35").Value = "L2100"
$ws.Range("D35").Value = "13-Dec-2025 1:34 PM"
$ws.Range("E35").Value = "FAIL"
$ws.Range("F35").Value = "1. TCP DL (Mbps)`n2. Network Detach time(ms)`n3. Network Detach Success rate`n4. Round trip time or Latency(ms)`n5. VOLTE setup SR`n6. Check functioning of ViLTE"
$ws.Range("G35").Value = "1. Static DL`n2. Static Cell Reselection`n3. Static Cell Reselection`n4. Static Ping`n5. Drive Volte SCMO`n6. Static Video call MO"
$ws.Range("H35").Value = "1. The peak value is not meeting the acceptance criteria. Please redo the static test in the site’s main lobe with good coverage, and use the Okhla Speedtest while running the script.`n2. The Detach Request count is higher than the Detach Accept count. The acceptance criteria require a 100 percent match. Please exclude this logfile, create a new one, and ensure that the Detach Request equals the Detach Accept in the Event tab.`n3. Exclude the current logfile and create a new one. Ensure that the Detach Request matches the Detach Accept. Verify this in the Event tab of AZQ before saving the logfile`n4. The Ping is above the acceptance criteria , Kinldy Exclud the logfile and Create new and verify the ping should bhe 40 to 50 ms , before saving the logfile`n5. Add a VoLTE short-call drive across all sectors. Do not disconnect the call manually, as it will be counted as a call drop. Each sector must have at least one successful MO session setup.`n6. Do not use WhatsApp for this test. Perform a manual VoLTE video call after running the script."

# --- Row 36: TALAK1_BLG_P47 ---
$ws.Range("A36").Value = "KK"
$ws.Range("B36").Value = "TALAK1_BLG_P47"
$ws.Range("C36").Value = "L2100"
$ws.Range("D36").Value = "15-Dec-2025 6:57 AM"
$ws.Range("E36").Value = "FAIL"
$ws.Range("F36").Value = "1. PCI`n2. RSRP (Average Value)`n3. SINR (Average Value)`n4. TCP DL (Mbps)`n5. TCP UL (Mbps)`n6. Network Detach time(ms)`n7. Network Detach Success rate`n8. Round trip time or Latency(ms)`n9. VOLTE setup SR`n10. Check functioning of ViLTE"
$ws.Range("G36").Value = "1. Static All`n2. Drive DL Ftp`n3. Drive DL Ftp`n4. Static DL`n5. Static UL`n6. Static Cell Reselection`n7. Static Cell Reselection`n8. Static Ping`n9. Drive Volte SCMO`n10. Static Video call MO"
$ws.Range("H36").Value = "1. Since the serving cell of the failed sector is different from the recorded PCI, the report is displaying the recorded PCI as the failed sector. Please exclude the current logfile and redo the test with the correct serving PCI.`n2. It appears that the PCI recorded in the site database does not match the actual serving PCI observed on the site, or the failed sector has no corresponding drive data. Request you to please verify the serving cells in the DL drive along with the PCI details in the site database.`n3. It appears that the PCI recorded in the site database does not match the actual serving PCI observed on the site, or the failed sector has no corresponding drive data. Request you to please verify the serving cells in the DL drive along with the PCI details in the site database.`n4. The peak value is not meeting the acceptance criteria. Please redo the static test in the site’s main lobe with good coverage, and use the Okhla Speedtest while running the script.`n5. The peak value is not meeting the acceptance criteria. Please redo the static test in the site’s main lobe with strong coverage. While performing the UL static test, start uploading a long video on YouTube to ensure proper uplink loading.`n6. The Detach Request count is higher than the Detach Accept count. The acceptance criteria require a 100 percent match. Please exclude this logfile, create a new one, and ensure that the Detach Request equals the Detach Accept in the Event tab.`n7. Exclude the current logfile and create a new one. Ensure that the Detach Request matches the Detach Accept. Verify this in the Event tab of AZQ before saving the logfile`n8. The Ping is above the acceptance criteria , Kinldy Exclud the logfile and Create new and verify the ping should bhe 40 to 50 ms , before saving the logfile`n9. Add a VoLTE short-call drive across all sectors. Do not disconnect the call manually, as it will be counted as a call drop. Each sector must have at least one successful MO session setup.`n10. Do not use WhatsApp for this test. Perform a manual VoLTE video call after running the script."

# --- Match the formatting (borders/font/alignment/wrap) used by the rows above ---
# by copying row 34's formats onto the two new rows (values were already set above).
$ws.Range("A34:H34").Copy()
$ws.Range("A35:H35").PasteSpecial(-4122)
$ws.Range("A34:H34").Copy()
$ws.Range("A36:H36").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row heights to fit the wrapped remark text ---
$ws.Rows.Item(35).RowHeight = 114
$ws.Rows.Item(36).RowHeight = 218.5

# --- Move the active selection down to the newly added row, like the saved workbook ---
$ws.Range("B35").Select()
